# Journal de travail: add three new rows (87, 88, 89) describing the work
# done around the comments feature, and grow the table / dimension /
# selection accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 87 --------------------------------------------------------------
# A87 already exists as an (empty) cell formatted as a date (style s="4"),
# so simply assigning a value to it keeps that date style.
$ws.Range("A87").Value() = 44705
$ws.Range("B87").Value() = "Réalisation"
$ws.Range("C87").Value() = 0.25
$ws.Range("D87").Value() = "Visite chef de projet -> Deboggage de l'ajout de commentaires"

# --- Row 88 ----------------------------------------------------------------
# A88 is a brand new cell - copy the date formatting from A87 before setting
# its value so it keeps the same (date, wrap-text) number format.
$ws.Range("A87").Copy($ws.Range("A88"))
$ws.Range("A88").Value() = 44705
$ws.Range("B88").Value() = "Réalisation"
$ws.Range("C88").Value() = 1
$ws.Range("D88").Value() = "Récupération des commentaires pour un article"

# --- Row 89 ----------------------------------------------------------------
$ws.Range("A87").Copy($ws.Range("A89"))
$ws.Range("A89").Value() = 44705
$ws.Range("B89").Value() = "Réalisation"
$ws.Range("C89").Value() = 0.5
$ws.Range("D89").Value() = "Édition de la doc selon modifications requises pour l'ajout de commentaires"

# Row 89 wraps onto two lines in the original file (ht="30")
$ws.Rows(89).RowHeight() = 30

# --- Grow the structured table (Tableau1) and its autofilter ---------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:F89"))

# --- Update the visible selection ------------------------------------------
$ws.Range("E89").Select()
